$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Reminder (days)" column -----------------------------------------
# Header for the new column F.
$ws.Range("F1").Value = "Reminder (days)"

# Reminder value (in days) for the existing event row.
$ws.Range("F3").Value = 7

# Widen column E (Description) and size the new column F to fit the
# "Reminder (days)" header text, matching the manual column-drag the user
# performed in Excel.
$ws.Columns.Item(5).ColumnWidth = 16.73
$ws.Columns.Item(6).ColumnWidth = 14.45

# --- Updated event ID -------------------------------------------------------
# The row's ID changed (e.g. the event was recreated with a new random ID).
# Force the cell to Text first so the 19-digit number round-trips exactly
# instead of being coerced into a double (which would lose precision).
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "3403722335677977069"

# --- Selection / active cell ------------------------------------------------
# The user's cursor ended up on K11 when the workbook was last saved.
$ws.Range("K11").Select()
